$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zone 1 repeaters (VK1RBM / "Z1" block) didn't have TX/RX frequencies filled
# in yet - add the 70cm simplex/duplex pairs for rows 2-4 (columns D = TX,
# E = RX), matching the rest of the table's layout.
$ws.Range("D2").Value = 439.125
$ws.Range("E2").Value = 439.125

$ws.Range("D3").Value = 431
$ws.Range("E3").Value = 431

$ws.Range("D4").Value = 439.175
$ws.Range("E4").Value = 439.175

# Start page numbering from 0 on print-out.
$ws.PageSetup.FirstPageNumber = 0

# Leave the selection on E5 (bottom-right pane of the frozen header/column)
# instead of the prior end-of-sheet selection.
[void]$ws.Range("E5").Select()
